# Generate Report for Handback
# - Updates the "Status" text on the Overview + per-locale sheets from
#   "Ready for handoff" to "Handed back: in sync with en-US".
# - Records that the handback round-trip is complete for a.md / b.md by
#   filling in "Latest Target File" (E) / "Latest Handback File" (F) with
#   the same source file + handoff xlf that were already linked, and
#   stamping a new "Latest Handback DateTime" (G).

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---- Overview sheet -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusHandedBack
$overview.Range("C2").Value = $statusHandedBack
$overview.Range("B3").Value = $statusHandedBack
$overview.Range("C3").Value = $statusHandedBack

# ---- Per-locale sheets ------------------------------------------------
$locales = @(
  @{ Sheet = "zh-cn"; Xlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8a2c0cc557cecf02021c8cd24e6a1d8586bf10ad/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; HandbackDate = "2016-02-06 03:53:33" }
  @{ Sheet = "de-de"; Xlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bfee945ab65e0090e74b21a64c988627087ad1f7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; HandbackDate = "2016-02-06 03:53:52" }
)

$aMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/867c1e5787c73d7d9e8f8320ea769e541091eed4/e2e/a.md"

foreach ($locale in $locales) {
  $ws = $wb.Worksheets.Item($locale.Sheet)

  foreach ($row in 2, 3) {
    # Status -> handed back
    $ws.Cells.Item($row, 2).Value = $statusHandedBack

    # E: Latest Target File = a.md (same source that was handed off)
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 5), $aMdUrl, "", "", "a.md")

    # F: Latest Handback File = the same xlf that was the handoff target
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 6), $locale.Xlf, "", "", $locale.Xlf)

    # G: Latest Handback DateTime
    $ws.Cells.Item($row, 7).Value = $locale.HandbackDate
  }
}
